# Add "Area" (column G) and "Atotal" (column H) to the discharge sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels (new shared strings "Area" / "Atotal")
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Per-segment cross-sectional area, mirroring the existing Q (discharge) columns:
#   G2  -> (D2-0)*B2/100          (first segment uses 0 as the lower bound)
#   G3  -> (D3-D2)*B3/100
#   G4:G11 -> (D{row}-D{row-1})*B{row}/100   (shared formula, like columns D/E)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

# Total cross-sectional area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Match the author's final selection (cell H2)
[void]$ws.Range("H2").Select()
